# calibrate probe sets for 2024 fieldwork
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header: "2024 experiments" ---
$ws.Range("A26").Value = "2024 experiments"

# --- Box5probes sub-table (rows 27-31) ---
$ws.Range("A27").Value = "Box5probes"
$ws.Range("B27").Value = "probe"
$ws.Range("C27").Value = "Shedd_coral_traceable"
$ws.Range("D27").Value = "offset"

$ws.Range("A28").Value = "T1"
$ws.Range("B28").Value = 22.9
$ws.Range("C28").Value = 23.04
$ws.Range("D28").Formula = "=B28-C28"

$ws.Range("A29").Value = "T2"
$ws.Range("B29").Value = 22.8
$ws.Range("C29").Value = 23.04

$ws.Range("A30").Value = "T3"
$ws.Range("B30").Value = 22.6
$ws.Range("C30").Value = 23.04

$ws.Range("A31").Value = "T4"
$ws.Range("B31").Value = 22.8
$ws.Range("C31").Value = 23.04

# Fill D29:D31 together so Excel records it as one shared formula group
$ws.Range("D29:D31").Formula = "=B29-C29"

# --- Box6probes sub-table (rows 33-37) ---
$ws.Range("A33").Value = "Box6probes"
$ws.Range("B33").Value = "probe"
$ws.Range("C33").Value = "Shedd_coral_traceable"
$ws.Range("D33").Value = "offset"

$ws.Range("A34").Value = "T1"
$ws.Range("B34").Value = 22.1
$ws.Range("C34").Value = 22.47
$ws.Range("D34").Formula = "=B34-C34"

$ws.Range("A35").Value = "T2"
$ws.Range("B35").Value = 21.8
$ws.Range("C35").Value = 22.47

$ws.Range("A36").Value = "T3"
$ws.Range("B36").Value = 22.2
$ws.Range("C36").Value = 22.47

$ws.Range("A37").Value = "T4"
$ws.Range("B37").Value = 22.1
$ws.Range("C37").Value = 22.47

# Fill D35:D37 together so Excel records it as one shared formula group
$ws.Range("D35:D37").Formula = "=B35-C35"

# --- Update selection to match the author's final cursor position ---
$ws.Range("E36").Select()
